# Apply the "stuff at the bottom of the sheets" commit:
#  1. Fill in the missing pair_kind ("generic") value for the four
#     practice header rows (J2:J5).
#  2. Append a new "stim details" block (rows 27-36) describing the
#     find-images / video / audio counts needed per month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the missing pair_kind values for the practice rows ---
$ws.Range("J2:J5").Value = "generic"

# --- 2. New "stim details" block ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# month / word_type pairs for rows 29-36
$stimRows = @(
  @(6, "video"),
  @(6, "video"),
  @(7, "video"),
  @(7, "video"),
  @(6, "audio"),
  @(6, "audio"),
  @(7, "audio"),
  @(7, "audio")
)

for ($i = 0; $i -lt $stimRows.Count; $i++) {
  $row = 29 + $i
  $ws.Cells.Item($row, 1).Value = $stimRows[$i][0]
  $ws.Cells.Item($row, 2).Value = $stimRows[$i][1]
}
